$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginPage")

# Register the built-in "Hyperlink" named cell style in the workbook
# (applied to a scratch cell, then cleared back to Normal).
$ws1.Range("Z1").Style = "Hyperlink"
$ws1.Range("Z1").Style = "Normal"
$ws1.Range("Z1").Clear()

$ws1.Range("C1").Value = "link"
$ws1.Range("C2").Value = "google"
$ws1.Range("C2").Font.ThemeColor = 11
